# Applies the "cat-gangster" content refresh described in the commit:
#   - New page title / meta title
#   - Reworded pros ("What we like") bullets
#   - Reworded cons ("What we don't like") bullets
#   - New bold title recap + new italic meta description
#
# NOTE: the six "What we like"/"What we don't like" bullet paragraphs each
# start with a stray *empty* run (<w:r/>) immediately before the run that
# carries the visible text, e.g.
#   <w:p><w:pPr><w:pStyle w:val="ListBullet"/>...</w:pPr><w:r/><w:r><w:t>OLD</w:t></w:r></w:p>
# A plain Find/Replace on paragraphs that carry a <w:pPr> rebuilds the
# paragraph's run list and silently drops that leading empty run (verified
# experimentally against this runtime). The title/bold-recap/italic
# paragraphs have no <w:pPr> and are unaffected, so they can use a plain
# Find & Replace. For the six bulleted paragraphs we instead locate the
# paragraph via Find and rewrite it with InsertXML, using the exact
# original <w:pPr>/run shape so only the visible text changes.

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Replace-ListItem($oldText, $newText) {
    # Paragraphs styled "ListBullet" with an empty leading run:
    #   <w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing .../><w:ind .../></w:pPr><w:r/><w:r><w:t>OLD</w:t></w:r></w:p>
    # Locate the paragraph by scanning $d.Paragraphs directly -- deriving a
    # Paragraph from an arbitrary Find-matched Range proved unreliable in
    # this runtime (it returned the wrong / stale paragraph), whereas
    # indexing $d.Paragraphs directly is accurate.
    $target = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$oldText*") {
            $target = $p
            break
        }
    }
    if ($null -eq $target) {
        throw "Text not found: $oldText"
    }
    $xml = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>$newText</w:t></w:r></w:p>"
    $target.Range.InsertXML($xml) | Out-Null
}

# 1) Main title (Heading1) + identical bold title recap near the end of the
#    document -- both paragraphs lack a leading empty run merge issue
#    (Heading1 has a single run; the bold recap has no <w:pPr>), and since
#    both runs carry identical old text, one ReplaceAll updates both.
Replace-Text "Play Cat Gangster Free Slot Game Review | High 5 Games" "Play Cat Gangster Free - Exciting Slot Game Review"

# 2) "What we like" bullets
Replace-ListItem "Multiway payline mechanism allows for more winning combinations" "Unique multiways payline mechanism"
Replace-ListItem "Two special symbols trigger bonus features" "Exciting bonus features with free spins"
Replace-ListItem "Players can choose their level of volatility to suit their playing style" "Customizable volatility options"
Replace-ListItem "Compatible with a range of mobile devices" "Visually appealing theme and graphics"

# 3) "What we don't like" bullets
Replace-ListItem "The biggest payout in the base game is relatively low" "Limited base game payout"
Replace-ListItem "Free spins bonus feature may not trigger frequently" "Not available on all platforms"

# 4) Italic meta description at the very end
Replace-Text "Read our review of Cat Gangster, a fun online slot game by High 5 Games. Play for free and discover the game's multiway payline mechanism and bonus features." "A detailed review of Cat Gangster slot game, including gameplay mechanics, bonus features, and potential for wins. Play for free!"

Write-Output "Done"
